$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 17, pushing the existing
# rows (old 17-89) down by one (new 18-90).
$ws.Rows("17:17").Insert()

# Fill in the new row 17 with the new record's data.
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value = 44749
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100102
$ws.Cells.Item(17, 8).Value = "Cítricos"
$ws.Cells.Item(17, 9).Value = 100102005
$ws.Cells.Item(17, 10).Value = "Naranja"
$ws.Cells.Item(17, 11).Value = "Fukumoto"
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 300
$ws.Cells.Item(17, 14).Value = 650
$ws.Cells.Item(17, 15).Value = 700
$ws.Cells.Item(17, 16).Value = 675
$ws.Cells.Item(17, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 675
$ws.Cells.Item(17, 20).Value = 1
